$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update stats for 2025-12 (row 25)
$ws.Range("B25").Value = 6482
$ws.Range("D25").Value = 6033091
$ws.Range("E25").Value = 930.7452946621413
$ws.Range("F25").Value = 10.03225258869462
$ws.Range("H25").Value = 26.3515304996224
